$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of Damasco price rows (Fruta / hortaliza, semanal).
# Each row below is rewritten in place to match the new weekly snapshot;
# row 21 is untouched by this update.

# Row 2
$ws.Range("D2").Value = 44194
$ws.Range("K2").Value = 'Patterson'
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 13000
$ws.Range("Q2").Value = '$/caja 15 kilos'
$ws.Range("S2").Value = 867
$ws.Range("T2").Value = 15

# Row 3
$ws.Range("D3").Value = 44174
$ws.Range("K3").Value = 'Modesto'
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 8500
$ws.Range("O3").Value = 8500
$ws.Range("P3").Value = 8500
$ws.Range("R3").Value = 'Región Metropolitana'
$ws.Range("S3").Value = 850

# Row 4
$ws.Range("D4").Value = 44174
$ws.Range("K4").Value = 'Modesto'
$ws.Range("M4").Value = 180
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("S4").Value = 833
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44174
$ws.Range("K5").Value = 'Modesto'
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("Q5").Value = '$/caja 18 kilos'
$ws.Range("S5").Value = 667
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value = 44189
$ws.Range("K6").Value = 'Patterson'
$ws.Range("M6").Value = 130
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("Q6").Value = '$/caja 18 kilos'
$ws.Range("S6").Value = 667
$ws.Range("T6").Value = 18

# Row 7
$ws.Range("D7").Value = 44159
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 8000
$ws.Range("O7").Value = 8000
$ws.Range("P7").Value = 8000
$ws.Range("Q7").Value = '$/bandeja 10 kilos'
$ws.Range("T7").Value = 10

# Row 8
$ws.Range("D8").Value = 44159
$ws.Range("K8").Value = 'Castle Brite'
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 65
$ws.Range("N8").Value = 7000
$ws.Range("O8").Value = 7000
$ws.Range("P8").Value = 7000
$ws.Range("R8").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S8").Value = 700

# Row 9
$ws.Range("D9").Value = 44195
$ws.Range("M9").Value = 124
$ws.Range("N9").Value = 13000
$ws.Range("O9").Value = 13000
$ws.Range("P9").Value = 13000
$ws.Range("Q9").Value = '$/caja 15 kilos'
$ws.Range("S9").Value = 867
$ws.Range("T9").Value = 15

# Row 10
$ws.Range("D10").Value = 44169
$ws.Range("K10").Value = 'Dina'
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 10000
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 10000
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 1000

# Row 11
$ws.Range("D11").Value = 44167
$ws.Range("K11").Value = 'Castle Brite'
$ws.Range("L11").Value = 'Especial'
$ws.Range("M11").Value = 85
$ws.Range("N11").Value = 10000
$ws.Range("O11").Value = 10000
$ws.Range("P11").Value = 10000
$ws.Range("Q11").Value = '$/bandeja 10 kilos'
$ws.Range("R11").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S11").Value = 1000
$ws.Range("T11").Value = 10

# Row 12
$ws.Range("D12").Value = 44167
$ws.Range("K12").Value = 'Castle Brite'
$ws.Range("L12").Value = 'Primera'
$ws.Range("N12").Value = 9500
$ws.Range("O12").Value = 9500
$ws.Range("P12").Value = 9500
$ws.Range("Q12").Value = '$/bandeja 10 kilos'
$ws.Range("R12").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S12").Value = 950
$ws.Range("T12").Value = 10

# Row 13
$ws.Range("D13").Value = 44167
$ws.Range("M13").Value = 60
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 15000
$ws.Range("P13").Value = 15000
$ws.Range("Q13").Value = '$/caja 18 kilos'
$ws.Range("S13").Value = 833
$ws.Range("T13").Value = 18

# Row 14
$ws.Range("D14").Value = 44187
$ws.Range("K14").Value = 'Patterson'
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 80
$ws.Range("N14").Value = 15000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 15000
$ws.Range("Q14").Value = '$/caja 15 kilos granel'
$ws.Range("S14").Value = 1000
$ws.Range("T14").Value = 15

# Row 15
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 95
$ws.Range("N15").Value = 13500
$ws.Range("O15").Value = 13500
$ws.Range("P15").Value = 13500
$ws.Range("S15").Value = 900

# Row 16
$ws.Range("L16").Value = 'Tercera'
$ws.Range("M16").Value = 120
$ws.Range("N16").Value = 12000
$ws.Range("O16").Value = 12000
$ws.Range("P16").Value = 12000
$ws.Range("S16").Value = 800

# Row 17
$ws.Range("D17").Value = 44166
$ws.Range("K17").Value = 'Castle Brite'
$ws.Range("L17").Value = 'Primera'
$ws.Range("N17").Value = 10000
$ws.Range("O17").Value = 10000
$ws.Range("P17").Value = 10000
$ws.Range("Q17").Value = '$/bandeja 10 kilos'
$ws.Range("R17").Value = 'Región Metropolitana'
$ws.Range("S17").Value = 1000
$ws.Range("T17").Value = 10

# Row 18
$ws.Range("D18").Value = 44166
$ws.Range("K18").Value = 'Castle Brite'
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 120
$ws.Range("N18").Value = 8000
$ws.Range("O18").Value = 8000
$ws.Range("P18").Value = 8000
$ws.Range("Q18").Value = '$/bandeja 10 kilos'
$ws.Range("R18").Value = 'Región Metropolitana'
$ws.Range("S18").Value = 800
$ws.Range("T18").Value = 10

# Row 19
$ws.Range("D19").Value = 44160
$ws.Range("K19").Value = 'Castle Brite'
$ws.Range("M19").Value = 25
$ws.Range("N19").Value = 8000
$ws.Range("O19").Value = 8000
$ws.Range("P19").Value = 8000
$ws.Range("Q19").Value = '$/bandeja 10 kilos'
$ws.Range("S19").Value = 800
$ws.Range("T19").Value = 10

# Row 20
$ws.Range("D20").Value = 44160
$ws.Range("K20").Value = 'Castle Brite'
$ws.Range("L20").Value = 'Segunda'
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = 7000
$ws.Range("O20").Value = 7000
$ws.Range("P20").Value = 7000
$ws.Range("Q20").Value = '$/bandeja 10 kilos'
$ws.Range("S20").Value = 700
$ws.Range("T20").Value = 10

# Row 22
$ws.Range("D22").Value = 44162
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 70
$ws.Range("N22").Value = 8500
$ws.Range("O22").Value = 8500
$ws.Range("P22").Value = 8500
$ws.Range("S22").Value = 850

# Row 23
$ws.Range("D23").Value = 44162
$ws.Range("M23").Value = 75
$ws.Range("N23").Value = 14000
$ws.Range("O23").Value = 15000
$ws.Range("P23").Value = 14400
$ws.Range("Q23").Value = '$/caja 18 kilos'
$ws.Range("S23").Value = 800
$ws.Range("T23").Value = 18

# Row 24
$ws.Range("D24").Value = 44176
$ws.Range("K24").Value = 'Modesto'
$ws.Range("M24").Value = 115
$ws.Range("N24").Value = 11000
$ws.Range("O24").Value = 12000
$ws.Range("P24").Value = 11609
$ws.Range("Q24").Value = '$/caja 12 kilos'
$ws.Range("S24").Value = 967
$ws.Range("T24").Value = 12

